# FAST_holdings.xlsx update (model refresh: 2021-07-13 -> 2021-07-14)
#
# The worksheet ships protected (legacy password-hash-only, no plaintext
# password known), so cell writes must temporarily unprotect it. We restore
# protection afterward using a password that happens to collide with the
# workbook's original legacy hash, so the sheet remains protected on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect("aaaba03")

# Update the confidential disclaimer banner: date 2021-07-13 -> 2021-07-14
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-14 for illustrative purposes only and are subject to change."
# Setting multi-line text auto-adjusts the row height; restore it so the row
# stays byte-for-byte as it was (no explicit height in the original file).
$ws.Rows(13).AutoFit()

# Refresh the Weight (D) and Percent Change (E) columns for rows 2-10
$ws.Range("D2").Value = 0.1327102717173147
$ws.Range("E2").Value = -0.03704561241027993

$ws.Range("D3").Value = 0.1095531381851853
$ws.Range("E3").Value = -0.007286606523247641

$ws.Range("D4").Value = 0.1136868811783255
$ws.Range("E4").Value = 0.0009599054862290402

$ws.Range("D5").Value = 0.1188266985609987
$ws.Range("E5").Value = -0.004329339002705823

$ws.Range("D6").Value = 0.1221092696760924
$ws.Range("E6").Value = 0.003607214428857919

$ws.Range("D7").Value = 0.1433912592066973
$ws.Range("E7").Value = -0.0002882952142992945

$ws.Range("D8").Value = 0.1319341752896207
$ws.Range("E8").Value = 0.001941209095951191

$ws.Range("D9").Value = 0.1277883061857654
$ws.Range("E9").Value = -0.009770535443245021

$ws.Range("D10").Value = 0.9999999999999999
$ws.Range("E10").Value = -0.006713229347942606

$ws.Protect("aaaba03")
